# The sheet is a weekly/biweekly price log for "Coliflor" at Femacal de
# La Calera. This edit adds one new sampling date (two rows: "Primera"
# and "Segunda" quality grades) at the top of the historical data block,
# pushing all the existing rows (318..379) down by two rows (to 320..381).
#
# Inserting at row 318 (rather than appending at the end) keeps the data
# block in its established most-recent-first order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current first data row of
# the historical block (row 318), shifting everything below it down.
$ws.Rows.Item(318).Insert()
$ws.Rows.Item(318).Insert()

# New row 318: "Primera" quality grade for the new date.
$ws.Cells.Item(318, 1).Value = 3
$ws.Cells.Item(318, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 44476
$ws.Cells.Item(318, 5).Value = 5
$ws.Cells.Item(318, 6).Value = 100112008
$ws.Cells.Item(318, 7).Value = "Coliflor"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 1350
$ws.Cells.Item(318, 11).Value = 650
$ws.Cells.Item(318, 12).Value = 650
$ws.Cells.Item(318, 13).Value = 650
$ws.Cells.Item(318, 14).Value = "`$/unidad"
$ws.Cells.Item(318, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(318, 16).Value = 650
$ws.Cells.Item(318, 17).Value = 1
$ws.Cells.Item(318, 18).Value = "Hortaliza"

# New row 319: "Segunda" quality grade for the same new date.
$ws.Cells.Item(319, 1).Value = 3
$ws.Cells.Item(319, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(319, 3).Value = "Coquimbo"
$ws.Cells.Item(319, 4).Value = 44476
$ws.Cells.Item(319, 5).Value = 5
$ws.Cells.Item(319, 6).Value = 100112008
$ws.Cells.Item(319, 7).Value = "Coliflor"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Segunda"
$ws.Cells.Item(319, 10).Value = 1300
$ws.Cells.Item(319, 11).Value = 550
$ws.Cells.Item(319, 12).Value = 550
$ws.Cells.Item(319, 13).Value = 550
$ws.Cells.Item(319, 14).Value = "`$/unidad"
$ws.Cells.Item(319, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(319, 16).Value = 550
$ws.Cells.Item(319, 17).Value = 1
$ws.Cells.Item(319, 18).Value = "Hortaliza"
